# Fix diagrams for Logic
# Repositions/resizes several label textboxes and connectors on slide 3
# (and one on slide 6), and right-aligns the text in a few of the
# relocated labels.
#
# NOTE on the literal numbers below: the host's Shape.Left/Top/Width/Height
# setters marshal through a single-precision (float32) point value and then
# floor() the EMU conversion, so a naive `emu / 12700` literal can land the
# stored EMU one unit short of the target. The constants used here were
# solved so that floor(float32(pt) * 12700) reproduces the exact target EMU
# from the target OOXML.

function Get-ShapeById {
    param($slide, [int]$id)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Shape 33 "TextBox 32" - post(TaskManagerChangedEvent) label
$shp33 = Get-ShapeById $s3 33
$shp33.Left = 446.7777252354331
$shp33.Top = 132.0
$shp33.Width = 240.0
$shp33.Height = 16.964095188188978
$shp33.TextFrame.TextRange.ParagraphFormat.Alignment = 3

# Shape 39 "Rectangle 62" - EventsCenter box
$shp39 = Get-ShapeById $s3 39
$shp39.Left = 618.0

# Shape 40 "Straight Connector 39"
$shp40 = Get-ShapeById $s3 40
$shp40.Left = 690.488372856693

# Shape 41 "Rectangle 40"
$shp41 = Get-ShapeById $s3 41
$shp41.Left = 684.8184509968504

# Shape 42 "Straight Arrow Connector 41"
$shp42 = Get-ShapeById $s3 42
$shp42.Left = 468.0
$shp42.Top = 156.0
$shp42.Width = 216.0

# Shape 44 "Straight Arrow Connector 43"
$shp44 = Get-ShapeById $s3 44
$shp44.Left = 468.0
$shp44.Top = 168.0
$shp44.Width = 212.81401834803148
$shp44.Height = 0.3455905511811024

# Shape 62 "TextBox 61" - post(TaskManagerChangedEvent) label (2nd occurrence)
$shp62 = Get-ShapeById $s3 62
$shp62.Left = 78.0
$shp62.Top = 378.0
$shp62.Width = 265.47291568582676
$shp62.TextFrame.TextRange.ParagraphFormat.Alignment = 3

# Shape 74 "TextBox 73" - handleTaskManagerChangedEvent() label
$shp74 = Get-ShapeById $s3 74
$shp74.Left = 372.0
$shp74.Width = 234.0
$shp74.TextFrame.TextRange.ParagraphFormat.Alignment = 3

# Shape 50 "TextBox 49" - handleTaskMangerChangedEvent() label
$shp50 = Get-ShapeById $s3 50
$shp50.Top = 424.83213808425194
$shp50.Width = 242.48220472440946

# ---------------------------------------------------------------------
# Slide 6
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)

# Shape 26 "TextBox 25" - execute("delete 1") label
$shp26 = Get-ShapeById $s6 26
$shp26.Left = 6.0
$shp26.Width = 126.92590721181102
